$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update mapped cell values to reflect the new PHIS (Person Health Information
# Search Results) mapping in place of the old PersonEvaluationSearchResults mapping.
$ws.Range("C2").Value = "/phisres-doc:PersonHealthInformationSearchResults/nc30:Person/phisres-ext:PersonPersistentIdentification"
$ws.Range("C3").Value = "/phisres-doc:PersonHealthInformationSearchResults/nc30:Person/phisres-ext:PersonTemporaryIdentification/nc30:IdentificationID"
$ws.Range("A4").Value = "Behavioral Health"
$ws.Range("C4").Value = "/phisres-doc:PersonHealthInformationSearchResults/phisres-ext:BehavioralHealthInformation/jxdm51:Evaluation/jxdm51:EvaluationDiagnosisDescriptionText"

# Rows 3 and 4 now need extra height to accommodate the longer wrapped text.
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30

# Column A widened slightly to fit "Behavioral Health".
$ws.Columns.Item(1).ColumnWidth = 14.83

# Update the active selection to C3, matching the author's last edit location.
$ws.Range("C3").Select()

# Update the saved window position of the workbook (the author dragged the
# Excel window to a different monitor/position before the last save).
$win = $excel.ActiveWindow
$win.Left = -28060
$win.Top = 2760
